$d = $word.ActiveDocument

# Locate the existing "1.04 ... debugging message." bullet item in the
# "Software Version" list so the new entry can be inserted right after it.
$rng = $d.Content
$found = $rng.Find.Execute("uSec debugging message.", $true, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the '1.04 ... debugging message.' paragraph."
}

$para = $rng.Paragraphs(1)

# Insert a new paragraph right after it; it inherits the same ListParagraph
# style / numbering (numId 3) from the paragraph it follows.
$para.Range.InsertParagraphAfter()

# Grab that freshly created paragraph and fill in its text.
$newPara = $para.Next()
$newPara.Range.Text = "1.05 Fixed the ‘F’ command that was broken in the last version."
